$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-29) holds the "Förändrad" date; bump each value by one day (45555 -> 45556)
$ws.Range("C2:C29").Value = 45556
